# CHS Yearly Financials update
# Inserts a new left-most data column (most recent period) before the
# existing column D, shifting the prior D:K data block to E:L, and fills
# the new column with the latest reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at D; everything from D:K moves to E:L.
$ws.Columns("D:D").Insert()

# 2) Seed the new column D with the same formatting/values the (now
#    shifted) column E carries, one block per table on the sheet, so that
#    number formats/styles (dates vs. integers) line up automatically.
$ws.Range("E7:E35").Copy($ws.Range("D7:D35"))
$ws.Range("E38:E77").Copy($ws.Range("D38:D77"))
$ws.Range("E80:E102").Copy($ws.Range("D80:D102"))

# Match the new column's width to its neighbour.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# 3) Overwrite the new column D with the actual new-period figures.

# -- Income Statement --
$ws.Range("D7").Value = 43498
$ws.Range("D8").Value = 2131100
$ws.Range("D9").Value = 1367700
$ws.Range("D10").Value = 763400
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 2087500
$ws.Range("D18").Value = 43700
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 135000
$ws.Range("D22").Value = 400
$ws.Range("D23").Value = 43300
$ws.Range("D24").Value = 12600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 30700
$ws.Range("D27").Value = 29800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 4900
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 34700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 34700

# -- Balance Sheet --
$ws.Range("D38").Value = 43498
$ws.Range("D41").Value = 124100
$ws.Range("D42").Value = 62000
$ws.Range("D43").Value = 21400
$ws.Range("D44").Value = 235200
$ws.Range("D45").Value = 42500
$ws.Range("D46").Value = 485200
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 370900
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 15200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1007000
$ws.Range("D57").Value = 143400
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 131800
$ws.Range("D60").Value = 275200
$ws.Range("D61").Value = 57500
$ws.Range("D62").Value = 94300
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 427100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 587100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 580000
$ws.Range("D77").Value = 0

# -- Cash Flow Statement --
$ws.Range("D80").Value = 43498
$ws.Range("D81").Value = 34700
$ws.Range("D83").Value = 91300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 158100
$ws.Range("D91").Value = -54200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -55900
$ws.Range("D96").Value = -43200
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -137700
$ws.Range("D101").Value = -500
$ws.Range("D102").Value = -35900
